# Natmi following Dr Hou advice
# Update the LR-pair statistics in rows 2-10 of the active sheet
# (columns E, G, H, I, J, K, M, N, O, P, Q, R, S, T) with recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  E=3; G=34.111822;         H=102.335466;       I=0.4228853893909983; J=0.4228853893909983; K=3; M=34.52052333333333; N=103.56157;  O=0.7684334662422598; P=0.7684334662422598; Q=1177.557947293513; R=10598.02152564162; S=0.3249592855929326;  T=0.3249592855929326}
    @{Row=3;  E=3; G=34.111822;         H=102.335466;       I=0.4228853893909983; J=0.4228853893909983; K=3; M=5.347618333333333; N=16.042855;  O=0.1190390091234806; P=0.1190390091234805; Q=182.4170047106033;  R=1641.75304239543;  S=0.05033985772590168; T=0.05033985772590167}
    @{Row=4;  E=3; G=34.111822;         H=102.335466;       I=0.4228853893909983; J=0.4228853893909983; K=3; M=5.055101333333334; N=15.165304;  O=0.1125275246342597; P=0.1125275246342597; Q=172.4387168746293;  R=1551.948451871664; S=0.04758624607216407; T=0.04758624607216407}
    @{Row=5;  E=3; G=34.88211266666666; H=104.646338;       I=0.4324347083490296; J=0.4324347083490295; K=3; M=34.52052333333333; N=103.56157;  O=0.7684334662422598; P=0.7684334662422598; Q=1204.148784225629; R=10837.33905803066; S=0.3322973018601055;  T=0.3322973018601055}
    @{Row=6;  E=3; G=34.88211266666666; H=104.646338;       I=0.4324347083490296; J=0.4324347083490295; K=3; M=5.347618333333333; N=16.042855;  O=0.1190390091234806; P=0.1190390091234805; Q=186.5362252016655;  R=1678.82602681499;  S=0.05147659919246979; T=0.05147659919246977}
    @{Row=7;  E=3; G=34.88211266666666; H=104.646338;       I=0.4324347083490296; J=0.4324347083490295; K=3; M=5.055101333333334; N=15.165304;  O=0.1125275246342597; P=0.1125275246342597; Q=176.3326142507502;  R=1586.993528256752; S=0.04866080729645433; T=0.04866080729645433}
    @{Row=8;  E=3; G=11.67052633333333; H=35.011579;        I=0.1446799022599722; J=0.1446799022599721; K=3; M=34.52052333333333; N=103.56157;  O=0.7684334662422598; P=0.7684334662422598; Q=402.8726766021144;  R=3625.854089419029; S=0.1111768787892218;  T=0.1111768787892217}
    @{Row=9;  E=3; G=11.67052633333333; H=35.011579;        I=0.1446799022599722; J=0.1446799022599721; K=3; M=5.347618333333333; N=16.042855;  O=0.1190390091234806; P=0.1190390091234805; Q=62.40952057978277;  R=561.6856852180449; S=0.0172225522051091;   T=0.01722255220510909}
    @{Row=10; E=3; G=11.67052633333333; H=35.011579;        I=0.1446799022599722; J=0.1446799022599721; K=3; M=5.055101333333334; N=15.165304;  O=0.1125275246342597; P=0.1125275246342597; Q=58.99569322833511;  R=530.961239055016;  S=0.0162804712656413;  T=0.0162804712656413}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 5).Value  = $item.E   # E: Ligand-expressing cells
    $ws.Cells.Item($r, 7).Value  = $item.G   # G: Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $item.H   # H: Ligand total expression value
    $ws.Cells.Item($r, 9).Value  = $item.I   # I: Ligand derived specificity of average expression value
    $ws.Cells.Item($r, 10).Value = $item.J   # J: Ligand derived specificity of total expression value
    $ws.Cells.Item($r, 11).Value = $item.K   # K: Receptor-expressing cells
    $ws.Cells.Item($r, 13).Value = $item.M   # M: Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $item.N   # N: Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $item.O   # O: Receptor derived specificity of average expression value
    $ws.Cells.Item($r, 16).Value = $item.P   # P: Receptor derived specificity of total expression value
    $ws.Cells.Item($r, 17).Value = $item.Q   # Q: Edge average expression weight
    $ws.Cells.Item($r, 18).Value = $item.R   # R: Edge total expression weight
    $ws.Cells.Item($r, 19).Value = $item.S   # S: Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value = $item.T   # T: Edge total expression derived specificity
}
